$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '81.144.98'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +6.19%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.214.59'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +4.26%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '212.70'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +7.08%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '629.85'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.35%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.282'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +30.40%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.999'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.05%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.592'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +7.14%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '3.213.65'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +4.24%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.598'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +31.52%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0000261'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +30.00%  '
$ws.Range('E13').Value = '  +1.83%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '3.799.02'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +3.84%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.30'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.79%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '32.26'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +9.48%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '80.721.70'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +5.79%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.173.98'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +3.10%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '14.42'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +6.79%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '3.04'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +10.64%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '9.31'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +2.18%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '443.60'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +14.80%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.22'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +15.27%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.97'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +8.11%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.335.14'
$ws.Range('D25').Style = 'Normal'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '76.69'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +5.70%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '4.73'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +2.75%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '11.00'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +5.96%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.999'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.16%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0000124'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +11.98%  '
$ws.Range('B31').Value = 'InternetComputer(DFINITY)'
$ws.Range('C31').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '9.05'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +7.85%  '
$ws.Range('B32').Value = 'Binance-PegBSC-USD'
$ws.Range('C32').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.997'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.42%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '569.30'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +12.64%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.49'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +3.07%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.152'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +14.57%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.00'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +3.38%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '23.24'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +11.04%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.125'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +21.20%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.00'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.15%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.408'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +7.75%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '20.81'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +3.69%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '164.70'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.60%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.75'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +9.98%  '
$ws.Range('B44').Value = 'Aave'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '192.38'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.71%  '
$ws.Range('B45').Value = 'USDe'
$ws.Range('C45').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.00'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.03%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.83'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +9.11%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.74'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +10.15%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.789'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.41%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.31'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +3.03%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '42.95'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +4.44%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '4.28'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +8.56%  '
